# Update send-email related Assets: add a new "CustomerName_SplitKeyword"
# asset (used to split the Egnyte file path to extract the customer name),
# and refresh the workbook view so the Assets sheet is the active tab.

$wb = $excel.ActiveWorkbook

$settings  = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")
$assets    = $wb.Worksheets.Item("Assets")

# --- Assets sheet: remove the blank spacer row before the email-support
# block (rows shift up by one) and append the new asset definition that
# used to be "lost" in that gap, right after Dispatcher_EmailBody_NoReport.
$assets.Rows.Item(28).Delete()

$assets.Cells.Item(32, 1).Value = "CustomerName_SplitKeyword"
$assets.Cells.Item(32, 2).Value = "CustomerName_SplitKeyword"
$assets.Cells.Item(32, 3).Value = "Contracts Concierge/Prod"
$assets.Cells.Item(32, 4).Value = "Keyword used to split the Egnyte file path to extract the customer name"

# --- Cosmetic row-height refresh (auto height recalculated by Excel for
# the wrapped description rows) on Settings / Constants.
$settings.Rows.Item(3).RowHeight = 43.5
$settings.Rows.Item(5).RowHeight = 29

$constants.Rows.Item(2).RowHeight = 29
$constants.Rows.Item(3).RowHeight = 43.5
$constants.Rows.Item(17).RowHeight = 29

# --- Window/view state: the Assets sheet is now the active tab, scrolled
# down near the newly-added asset row.
$assets.Activate()
$assets.Range("C35").Select()
